{"js": "// 1) Remove the \"Meta description\" paragraph (the 2nd paragraph of the\n//    document, right after the Heading1 title paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[1].delete();\nawait context.sync();\n\n// 2) Re-load paragraphs (indices shifted after the delete) and update the\n//    very last paragraph: insert a new bold paragraph right before it\n//    (\"Play FAUST Slot Game for Free | Game Review\") and replace the last\n//    paragraph's text with the meta-description sentence (keeping its\n//    existing italic formatting).\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst lastIndex = paragraphs2.items.length - 1;\nconst secondToLast = paragraphs2.items[lastIndex - 1];\nconst lastParagraph = paragraphs2.items[lastIndex];\n\n// Insert the new heading-like paragraph after the paragraph that currently\n// precedes the last one, so it lands directly before the last paragraph\n// without inheriting the last paragraph's italic run formatting.\nconst newParagraph = secondToLast.insertParagraph(\n  \"Play FAUST Slot Game for Free | Game Review\",\n  Word.InsertLocation.after\n);\nnewParagraph.style = \"Normal\";\nnewParagraph.font.bold = true;\n\n// Replace the text of the final paragraph in place (preserves its italic\n// run formatting and leading empty run).\nlastParagraph.insertText(\n  \"Read our review of the FAUST slot game. Play for free and enjoy its storybook design and expanding symbol free spin bonuses.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"Meta description\" paragraph (the 2nd paragraph of the\n#    document, right after the Heading1 title paragraph).\n$metaPara = $d.Paragraphs(2)\n$metaPara.Range.Delete()\n\n# 2) Insert a new bold paragraph (\"Play FAUST Slot Game for Free | Game\n#    Review\") right before the very last paragraph of the document, without\n#    disturbing the last paragraph's existing (italic) run formatting.\n$count = $d.Paragraphs.Count\n$secondToLast = $d.Paragraphs($count - 1)\n$secondToLast.Range.InsertParagraphAfter()\n\n$count2 = $d.Paragraphs.Count\n$newPara = $d.Paragraphs($count2 - 1)\n$newPara.Style = $d.Styles(\"Normal\")\n$newPara.Range.Text = \"Play FAUST Slot Game for Free | Game Review\"\n\n# Bold only the literal text (exclude the trailing paragraph mark) so no\n# stray formatting gets attached to the paragraph mark itself.\n$textRng = $newPara.Range\n$textRng.MoveEnd(1, -1)\n$textRng.Font.Bold = $true\n\n# 3) Replace the text of the final paragraph (the old AI image prompt) with\n#    the meta-description sentence, keeping its italic run formatting intact.\n$count3 = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($count3)\n$find = $lastPara.Range.Find\n$find.Text = \"Create a cartoon-style image featuring a happy Maya warrior with glasses, holding a potion and standing next to Faust. They are both standing in an alchemist's lab, surrounded by bubbling flasks and beakers. In the background, we see the silhouette of Mephistopheles lurking in the shadows. The image should convey a sense of excitement and adventure, with a touch of humor. The colors should be bold and vibrant, drawing the player's attention to the game.\"\n$find.Replacement.Text = \"Read our review of the FAUST slot game. Play for free and enjoy its storybook design and expanding symbol free spin bonuses.\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
